# Insert a new data row at row 99 (pushes existing rows 99-154 down to
# 100-155) on the "Fruta, Macroferia Regional de Talca - Mango" sheet,
# and fill the newly inserted row with a fresh weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 99; formatting (incl. the date
# number format on column D) is inherited from the surrounding rows.
$ws.Rows.Item(99).Insert()

$ws.Cells.Item(99, 1).Value = 5
$ws.Cells.Item(99, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(99, 3).Value = "Maule"
$ws.Cells.Item(99, 4).Value = 44873
$ws.Cells.Item(99, 5).Value = 7
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100108
$ws.Cells.Item(99, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(99, 9).Value = 100108002
$ws.Cells.Item(99, 10).Value = "Mango"
$ws.Cells.Item(99, 11).Value = "Sin especificar"
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 248
$ws.Cells.Item(99, 14).Value = 8000
$ws.Cells.Item(99, 15).Value = 9000
$ws.Cells.Item(99, 16).Value = 8806
$ws.Cells.Item(99, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(99, 18).Value = "Brasil"
$ws.Cells.Item(99, 19).Value = 2202
$ws.Cells.Item(99, 20).Value = 4
